$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.871.58"
$ws.Range("D3").Value = "2.349.19"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.668"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.69"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").Value = "2.702.88"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.904"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "2.350.52"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "43.811.11"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "78.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "176.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.128"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0748"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.92%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +18.31%  "
$ws.Range("E41").Value = "  -3.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +16.14%  "
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.200"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.83%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.105"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.83%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.39%  "
